# Daily attendance processing - 2025-11-10 12:41:05
#
# Column G ("Recorded By") holds comma-separated lists of recorder
# identities (e.g. "dnasr281@gmail.com, System"). This pass normalizes
# the ordering of those lists onto a fixed canonical priority order:
#   backup@backdoor.com < System < system < admin@admin.com < dnasr281@gmail.com
#
# Any cell whose comma-separated values are not already in that order
# gets rewritten with the values reordered accordingly. Single-value
# cells (and cells already in the correct order) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RecorderPriority($s) {
    if ($s.Equals("backup@backdoor.com")) { return 0 }
    elseif ($s.Equals("System")) { return 1 }
    elseif ($s.Equals("system")) { return 2 }
    elseif ($s.Equals("admin@admin.com")) { return 3 }
    elseif ($s.Equals("dnasr281@gmail.com")) { return 4 }
    else { return 99 }
}

function Sort-Recorders($val) {
    $parts = @($val -split ",\s*")
    $n = $parts.Count
    # Stable bubble sort by canonical priority (case-sensitive tokens).
    for ($i = 0; $i -lt $n; $i++) {
        for ($j = 0; $j -lt ($n - $i - 1); $j++) {
            $p1 = Get-RecorderPriority($parts[$j].Trim())
            $p2 = Get-RecorderPriority($parts[$j + 1].Trim())
            if ($p1 -gt $p2) {
                $tmp = $parts[$j]
                $parts[$j] = $parts[$j + 1]
                $parts[$j + 1] = $tmp
            }
        }
    }
    return [string]::Join(", ", $parts)
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G: Recorded By
    $val = $cell.Value2

    if ($null -eq $val -or $val -eq "") {
        continue
    }

    $parts = @($val -split ",\s*")
    if ($parts.Count -lt 2) {
        continue
    }

    $newVal = Sort-Recorders $val

    if (-not $newVal.Equals($val)) {
        $cell.Value2 = $newVal
    }
}
